# edit.ps1 - apply the "change structure in the PPt" commit to the deck.
#
# Summary of changes (see XML diff in task description):
#   1. Two presentation-level drawing guides are nudged to new positions.
#   2. Slide 6, shape "內容版面配置區 2" (id 4, the big left-hand file-tree
#      textbox) is resized/repositioned and gets two new bullet lines
#      inserted into its file-tree listing (after the 4th "index.html"
#      entry, right before the "Main webpage" index.html line):
#         ├── Data Structure Visualization.pptx         # Visualization slides
#         ├── favicon.ico                               # Website favicon
#   3. Slide 6, shape "內容版面配置區 2" (id 9, the small textbox on the
#      right) is nudged to a new position (size unchanged).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

# ---------------------------------------------------------------------
# 1) Drawing guides (presentation-level). PowerPoint's Guide.Position is
#    expressed in points; the OOXML p15:guide@pos attribute stores
#    eighths of a point, e.g. XML pos 2904 == 363pt, 2880 == 360pt.
#    Guarded with try/catch so the rest of the script still applies even
#    if this host/version doesn't surface the Guides collection.
# ---------------------------------------------------------------------
try {
    $guides = $p.Guides
    if ($guides -ne $null) {
        for ($i = 1; $i -le $guides.Count; $i++) {
            $guide = $guides.Item($i)
            if ($guide.Orientation -eq 1 -and [Math]::Round($guide.Position) -eq 363) {
                $guide.Position = 360
            } elseif ($guide.Orientation -eq 2 -and [Math]::Round($guide.Position) -eq 480) {
                $guide.Position = 60
            }
        }
    }
} catch {
    # Guides collection not supported by this host - skip silently.
}

# ---------------------------------------------------------------------
# 2) Resize/reposition the big file-tree textbox (shape id 4) and add
#    two new lines to its text.
# ---------------------------------------------------------------------
$shpTree = $s.Shapes.Item(1)

$shpTree.Left   = 349.7411193847656
$shpTree.Top    = 21.416458129882812
$shpTree.Width  = 396.7500915527344
$shpTree.Height = 507.3108825683594

$tr = $shpTree.TextFrame.TextRange
# Paragraph 35 (1-based) is "|            |-- index.html" for the tree/
# folder - the last of the four repeated index.html lines, right before
# the "Main webpage" index.html line (paragraph 36).
$anchorPara = $tr.Paragraphs(35)
$null = $anchorPara.InsertAfter("`r├── Data Structure Visualization.pptx         # Visualization slides`r├── favicon.ico                               # Website favicon")

# ---------------------------------------------------------------------
# 3) Reposition the small textbox (shape id 9) - size is unchanged.
# ---------------------------------------------------------------------
$shpSmall = $s.Shapes.Item(3)
$shpSmall.Left = 663.7616577148438
$shpSmall.Top  = 187.29551696777344

Write-Output "edit applied"
